$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1864951768488746
$ws.Range("C2").Value = 0.6045016077170418
$ws.Range("J2").Value = 0.009646302250803859
$ws.Range("P2").Value = 0.1446945337620579
$ws.Range("S2").Value = 0.05466237942122187
$ws.Range("B3").Value = 0.01052631578947368
$ws.Range("C3").Value = 0.005263157894736842
$ws.Range("J3").Value = 0.005263157894736842
$ws.Range("P3").Value = 0.7578947368421053
$ws.Range("S3").Value = 0.2210526315789474
$ws.Range("J4").Value = 0.01785714285714286
$ws.Range("P4").Value = 0.8035714285714286
$ws.Range("S4").Value = 0.1785714285714286
$ws.Range("B6").Value = 0.08294930875576037
$ws.Range("D6").Value = 0.02304147465437788
$ws.Range("F6").Value = 0.04608294930875576
$ws.Range("J6").Value = 0.1797235023041475
$ws.Range("O6").Value = 0.02304147465437788
$ws.Range("Q6").Value = 0.1935483870967742
$ws.Range("R6").Value = 0.08755760368663594
$ws.Range("S6").Value = 0.3640552995391705
$ws.Range("B7").Value = 0.08372093023255814
$ws.Range("D7").Value = 0.02790697674418605
$ws.Range("F7").Value = 0.05581395348837209
$ws.Range("J7").Value = 0.07441860465116279
$ws.Range("O7").Value = 0.009302325581395349
$ws.Range("Q7").Value = 0.1906976744186047
$ws.Range("R7").Value = 0.07441860465116279
$ws.Range("S7").Value = 0.4837209302325581
$ws.Range("B8").Value = 0.06951871657754011
$ws.Range("D8").Value = 0.0231729055258467
$ws.Range("E8").Value = 0.0017825311942959
$ws.Range("F8").Value = 0.0677361853832442
$ws.Range("J8").Value = 0.1016042780748663
$ws.Range("O8").Value = 0.0071301247771836
$ws.Range("Q8").Value = 0.1764705882352941
$ws.Range("R8").Value = 0.09090909090909091
$ws.Range("S8").Value = 0.4616755793226381
$ws.Range("B9").Value = 0.09278350515463918
$ws.Range("D9").Value = 0.02061855670103093
$ws.Range("F9").Value = 0.05670103092783505
$ws.Range("J9").Value = 0.08247422680412371
$ws.Range("O9").Value = 0.0154639175257732
$ws.Range("Q9").Value = 0.1804123711340206
$ws.Range("R9").Value = 0.1030927835051546
$ws.Range("S9").Value = 0.4484536082474227
$ws.Range("B10").Value = 0.1187214611872146
$ws.Range("D10").Value = 0.02207001522070015
$ws.Range("F10").Value = 0.0639269406392694
$ws.Range("J10").Value = 0.09208523592085235
$ws.Range("O10").Value = 0.0136986301369863
$ws.Range("Q10").Value = 0.2168949771689498
$ws.Range("R10").Value = 0.08143074581430745
$ws.Range("S10").Value = 0.3911719939117199
$ws.Range("G11").Value = 0.1349862258953168
$ws.Range("J11").Value = 0.1046831955922865
$ws.Range("K11").Value = 0.2148760330578512
$ws.Range("L11").Value = 0.5371900826446281
$ws.Range("S11").Value = 0.008264462809917356
$ws.Range("G12").Value = 0.72
$ws.Range("J12").Value = 0.24
$ws.Range("K12").Value = 0.005
$ws.Range("L12").Value = 0.03
$ws.Range("S12").Value = 0.005
$ws.Range("G13").Value = 0.58
$ws.Range("J13").Value = 0.36
$ws.Range("S13").Value = 0.06
$ws.Range("F15").Value = 0.005
$ws.Range("H15").Value = 0.165
$ws.Range("I15").Value = 0.065
$ws.Range("J15").Value = 0.39
$ws.Range("K15").Value = 0.075
$ws.Range("M15").Value = 0.03
$ws.Range("O15").Value = 0.015
$ws.Range("S15").Value = 0.255
$ws.Range("F16").Value = 0.02164502164502164
$ws.Range("H16").Value = 0.1991341991341991
$ws.Range("I16").Value = 0.08225108225108226
$ws.Range("J16").Value = 0.3506493506493507
$ws.Range("K16").Value = 0.1298701298701299
$ws.Range("M16").Value = 0.03463203463203463
$ws.Range("O16").Value = 0.0303030303030303
$ws.Range("S16").Value = 0.1515151515151515
$ws.Range("F17").Value = 0.01214574898785425
$ws.Range("H17").Value = 0.1882591093117409
$ws.Range("I17").Value = 0.08097165991902834
$ws.Range("J17").Value = 0.4048582995951417
$ws.Range("K17").Value = 0.1214574898785425
$ws.Range("M17").Value = 0.02024291497975709
$ws.Range("N17").Value = 0.002024291497975709
$ws.Range("O17").Value = 0.05060728744939271
$ws.Range("S17").Value = 0.1194331983805668
$ws.Range("F18").Value = 0.004651162790697674
$ws.Range("H18").Value = 0.213953488372093
$ws.Range("I18").Value = 0.06511627906976744
$ws.Range("J18").Value = 0.4186046511627907
$ws.Range("K18").Value = 0.09767441860465116
$ws.Range("M18").Value = 0.02325581395348837
$ws.Range("O18").Value = 0.09302325581395349
$ws.Range("S18").Value = 0.08372093023255814
$ws.Range("F19").Value = 0.01024590163934426
$ws.Range("H19").Value = 0.2342896174863388
$ws.Range("I19").Value = 0.07581967213114754
$ws.Range("J19").Value = 0.3545081967213115
$ws.Range("K19").Value = 0.1086065573770492
$ws.Range("M19").Value = 0.0157103825136612
$ws.Range("O19").Value = 0.06284153005464481
$ws.Range("S19").Value = 0.1379781420765027
